# EasyController BOM update — add heatsink/thermal paste, replace Arduino Nano
# with Teensy 3.2 (+5V linear reg, 14-pin headers), fix several designators /
# quantities / prices, and drop the now-unused 680uF cap line.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Part number / designator / qty / price / description columns (A:E) ---
# Row 2: Gate driver
$ws.Range("A2").Value2 = "497-6219-5-ND"
$ws.Range("B2").Value2 = "IC2,IC3,IC4"
$ws.Range("C2").Value2 = 3
$ws.Range("D2").Value2 = 1.66
$ws.Range("E2").Value2 = "Gate driver"

# Row 3: 1uF cap
$ws.Range("A3").Value2 = "445-173257-1-ND"
$ws.Range("B3").Value2 = "C1,C2,C3,C11,C12,C13,C14"
$ws.Range("C3").Value2 = 7
$ws.Range("D3").Value2 = 0.33
$ws.Range("E3").Value2 = "1uF cap"

# Row 4: NFET
$ws.Range("A4").Value2 = "IRFB7730PBF-ND"
$ws.Range("B4").Value2 = "Q1,Q2,Q3,Q4,Q5,Q6"
$ws.Range("C4").Value2 = 6
$ws.Range("D4").Value2 = 3.51
$ws.Range("E4").Value2 = "NFET"

# Row 5: TVS Diode
$ws.Range("A5").Value2 = "P6KE68CALFCT-ND"
$ws.Range("B5").Value2 = "D1"
$ws.Range("C5").Value2 = 1
$ws.Range("D5").Value2 = 0.47
$ws.Range("E5").Value2 = "TVS Diode"

# Row 6: 470uF cap (replaces 680uF cap)
$ws.Range("A6").Value2 = "565-3994-ND"
$ws.Range("B6").Value2 = "C4,C5,C10"
$ws.Range("C6").Value2 = 3
$ws.Range("D6").Value2 = 1.15
$ws.Range("E6").Value2 = "470uF cap"

# Row 7: 2.2k resistor
$ws.Range("A7").Value2 = "CF14JT2K20CT-ND"
$ws.Range("B7").Value2 = "R1,R2,R3,R4,R5,R6"
$ws.Range("C7").Value2 = 6
$ws.Range("D7").Value2 = 0.1
$ws.Range("E7").Value2 = "2.2k resistor"

# Row 8: 47nF cap
$ws.Range("A8").Value2 = "BC1082CT-ND"
$ws.Range("B8").Value2 = "C6,C7,C8,C9"
$ws.Range("C8").Value2 = 3
$ws.Range("D8").Value2 = 0.16
$ws.Range("E8").Value2 = "47nF cap"

# Row 9: Fuse holder
$ws.Range("A9").Value2 = "36-3557-2-ND"
$ws.Range("B9").Value2 = "J3"
$ws.Range("C9").Value2 = 1
$ws.Range("D9").Value2 = 1.03
$ws.Range("E9").Value2 = "Fuse holder"

# Row 10: Voltage Regulator
$ws.Range("A10").Value2 = "1470-4533-ND"
$ws.Range("B10").Value2 = "IC1"
$ws.Range("C10").Value2 = 1
$ws.Range("D10").Value2 = 8
$ws.Range("E10").Value2 = "Voltage Regulator"

# Row 11: 47k resistor
$ws.Range("A11").Value2 = "CF14JT47K0CT-ND"
$ws.Range("B11").Value2 = "R7,R8"
$ws.Range("C11").Value2 = 1
$ws.Range("D11").Value2 = 0.1
$ws.Range("E11").Value2 = "47k resistor"

# Row 12: Teensy 3.2 (replaces Arduino Nano)
$ws.Range("A12").Value2 = "1568-1231-ND"
$ws.Range("B12").Value2 = "MCU1"
$ws.Range("C12").Value2 = 1
$ws.Range("D12").Value2 = 22.5
$ws.Range("E12").Value2 = "Teensy 3.2"

# Row 13: 22 ohm resistor (replaces the old "15 pin female header" line)
$ws.Range("A13").Value2 = "CF14JT22R0CT-ND"
$ws.Range("B13").Value2 = "R9,R10,R11,R12,R13,R14"
$ws.Range("C13").Value2 = 6
$ws.Range("D13").Value2 = 0.1
$ws.Range("E13").Value2 = "22 ohm"

# Row 14: 5V Linear reg (new)
$ws.Range("A14").Value2 = "497-1184-1-ND"
$ws.Range("B14").Value2 = "IC5"
$ws.Range("C14").Value2 = 1
$ws.Range("D14").Value2 = 0.35
$ws.Range("E14").Value2 = "5V Linear reg"

# Row 15: 14 pin male header (new)
$ws.Range("A15").Value2 = "S1012EC-14-ND"
$null = $ws.Range("B15").ClearContents()
$ws.Range("C15").Value2 = 2
$ws.Range("D15").Value2 = 0.31
$ws.Range("E15").Value2 = "14 pin male header"

# Row 16: 14 pin female header (new)
$ws.Range("A16").Value2 = "S7012-ND"
$null = $ws.Range("B16").ClearContents()
$ws.Range("C16").Value2 = 2
$ws.Range("D16").Value2 = 0.91
$ws.Range("E16").Value2 = "14 pin female header"

# Row 17: stays blank in A:H, only the extended-price formula remains
$null = $ws.Range("A17:H17").ClearContents()

# Row 18 previously held the trailing shared-formula cell only; it goes away
# entirely now that the formula range ends at row 17.
$null = $ws.Range("A18:K18").ClearContents()

# --- Extended price column (I) and running total (K2) ---
$ws.Range("I2").Formula = "=D2*C2"
$ws.Range("I3").Formula = "=D3*C3"
$ws.Range("I4").Formula = "=D4*C4"
$ws.Range("I5").Formula = "=D5*C5"
$ws.Range("I6").Formula = "=D6*C6"
$ws.Range("I7").Formula = "=D7*C7"
$ws.Range("I8").Formula = "=D8*C8"
$ws.Range("I9").Formula = "=D9*C9"
$ws.Range("I10").Formula = "=D10*C10"
$ws.Range("I11").Formula = "=D11*C11"
$ws.Range("I12").Formula = "=D12*C12"
$ws.Range("I13").Formula = "=D13*C13"
$ws.Range("I14").Formula = "=D14*C14"
$ws.Range("I15").Formula = "=D15*C15"
$ws.Range("I16").Formula = "=D16*C16"
$ws.Range("I17").Formula = "=D17*C17"

$ws.Range("K2").Formula = "=SUM(I2:I99)"

# --- New rows 19/20: heatsink + thermal paste ---
$ws.Range("A19").Value2 = "345-1092-ND" + [char]0x200E + " "
$ws.Range("E19").Value2 = "Heatsink"

$ws.Range("A20").Value2 = [char]0x200E + "1168-1380-ND" + [char]0x200E + " "
$ws.Range("E20").Value2 = "Thermal paste"

# --- Selection cursor, matching the saved workbook state ---
$null = $ws.Range("G18").Select()
